$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.950.00"
$ws.Range("E2").Value = "  -2.26%  "

# Row 3
$ws.Range("D3").Value = "3.567.95"
$ws.Range("E3").Value = "  -3.05%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.19"
$ws.Range("E5").Value = "  -6.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.57"
$ws.Range("E6").Value = "  -3.24%  "

# Row 7
$ws.Range("D7").Value = "3.564.52"
$ws.Range("E7").Value = "  -3.08%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -1.90%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11
$ws.Range("E11").Value = "  -2.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("E12").Value = "  -1.66%  "

# Row 13
$ws.Range("E13").Value = "  -2.92%  "

# Row 14
$ws.Range("D14").Value = "4.166.08"
$ws.Range("E14").Value = "  -3.23%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.31"
$ws.Range("E15").Value = "  -1.26%  "

# Row 16
$ws.Range("D16").Value = "3.566.84"
$ws.Range("E16").Value = "  -3.03%  "

# Row 17
$ws.Range("D17").Value = "67.949.09"
$ws.Range("E17").Value = "  -2.27%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.117"
$ws.Range("E18").Value = "  -0.82%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  -0.21%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.72"
$ws.Range("E20").Value = "  -1.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "455.61"
$ws.Range("E21").Value = "  -2.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.62"
$ws.Range("E22").Value = "  -1.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.644"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.76"
$ws.Range("E24").Value = "  -2.55%  "

# Row 25
$ws.Range("D25").Value = "3.705.85"
$ws.Range("E25").Value = "  -3.19%  "

# Row 26
$ws.Range("E26").Value = "  +0.16%  "

# Row 27 - now PEPE (was InternetComputer(DFINITY))
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  -6.43%  "

# Row 28 - now InternetComputer(DFINITY) (was PEPE)
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.64"
$ws.Range("E28").Value = "  -2.52%  "

# Row 29
$ws.Range("E29").Value = "  -6.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.57"
$ws.Range("E30").Value = "  -3.33%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -2.69%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.05"
$ws.Range("E33").Value = "  -2.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.92"
$ws.Range("E34").Value = "  -3.88%  "

# Row 35 - now Kaspa (was NEARProtocol)
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.160"
$ws.Range("E35").Value = "  -1.89%  "

# Row 36 - now NEARProtocol (was Kaspa)
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.25"
$ws.Range("E36").Value = "  -3.24%  "

# Row 37
$ws.Range("D37").Value = "3.564.90"
$ws.Range("E37").Value = "  -2.93%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.12"
$ws.Range("E38").Value = "  -3.63%  "

# Row 39
$ws.Range("E39").Value = "  +0.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "177.35"
$ws.Range("E41").Value = "  -1.00%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0890"
$ws.Range("E42").Value = "  -1.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.64"
$ws.Range("E43").Value = "  -6.75%  "

# Row 44
$ws.Range("E44").Value = "  -5.11%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.893"
$ws.Range("E45").Value = "  -4.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.27"
$ws.Range("E46").Value = "  +6.86%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.17"
$ws.Range("E47").Value = "  -1.61%  "

# Row 48
$ws.Range("E48").Value = "  -4.02%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.72"
$ws.Range("E49").Value = "  -1.26%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.21"
$ws.Range("E50").Value = "  -5.42%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.02"
$ws.Range("E51").Value = "  -4.01%  "

